$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: new "Bit Manupulation" / Scaler question (magic numbers)
$ws.Range("A16").Value = "Bit Manupulation"
$ws.Range("B16").Value = "Scaler"
$ws.Range("C16").Value = "Given an integer A, find and return the Ath magic number.`nA magic number is defined as a number that can be expressed as a power of 5 or a sum of unique powers of 5.`nFirst few magic numbers are 5, 25, 30(5 + 25), 125, 130(125 + 5), $([char]0x2026)."
$ws.Range("D16").Value = "com.sumeet.dsa.array.Solution6"

# Row 17: new "Maths" / Scaler question (distinct primes sum)
$ws.Range("A17").Value = "Maths"
$ws.Range("B17").Value = "Scaler"
$ws.Range("C17").Value = "Given an even no , find if sum of any two distinct prime is equal to that no"
$ws.Range("D17").Value = "com.sumeet.dsa.array.Solution7"

# Match the existing "Question" column formatting (wrap text) used by the other rows
$ws.Range("C16:C17").WrapText = $true
$ws.Rows.Item(16).RowHeight = 90

# Reflect the author's final cursor position/selection
$ws.Range("D24").Select() | Out-Null
